# Vancouver Home Defensive Actions - cleaned defensive actions data
#
# The original sheet was exported from pandas with a 2-row MultiIndex
# header (row 1 = grouped/"Unnamed" labels with merged cells, row 2 =
# the real column names) and a blank row 3 separating the header from
# the data. This edit flattens the header into a single visible row,
# hides the old (now redundant) second header row and the blank
# separator row, de-duplicates the "Challenges"/"Tackles" labels down
# to short codes, fills in the previously-blank Tkl% cells with 0, and
# hides the trailing summary row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Flatten the header into row 1 -----------------------------------
# Remove the merged cells that grouped H1:L1 ("Tackles"), M1:P1
# ("Challenges") and Q1:S1 ("Blocks") so each column header can stand on
# its own again.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# --- Hide the now-redundant original column-name row and the blank
#     separator row beneath it -----------------------------------------
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true

# --- Fill in the previously-blank Tkl% cells with an explicit 0 -------
$ws.Range("O4").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("O18").Value = 0

# --- Hide the trailing "16 Players" summary row ------------------------
$ws.Rows.Item(20).Hidden = $true
